# Auto-generated Excel COM-interop script
# Applies updated market-price / profit values across the Leve profit sheets
# (columns H:N) to match the refreshed source data.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 5: H5,I5,J5,K5,L5,M5,N5
$ws.Range("H5").Value = 68.59999999999999
$ws.Range("I5").Value = 62.75
$ws.Range("J5").Value = 92
$ws.Range("K5").Value = 62.75
$ws.Range("L5").Value = 92
$ws.Range("M5").Value = 52.25
$ws.Range("N5").Value = -322
# Row 12: H12,I12,J12,K12,L12,M12,N12
$ws.Range("H12").Value = 473.22223
$ws.Range("I12").Value = 499.75
$ws.Range("J12").Value = 452
$ws.Range("K12").Value = 499.75
$ws.Range("L12").Value = 452
$ws.Range("M12").Value = -329.75
$ws.Range("N12").Value = -792
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 5408320
$ws.Range("I132").Value = 6669379
$ws.Range("J132").Value = 3781.1428
$ws.Range("K132").Value = 20008137
$ws.Range("L132").Value = 11343.4284
$ws.Range("M132").Value = -20005607
$ws.Range("N132").Value = -16403.4284
# Row 141: H141,I141,J141,K141,L141,M141,N141
$ws.Range("H141").Value = 758551.4
$ws.Range("I141").Value = 2475.6667
$ws.Range("J141").Value = 1325608.1
$ws.Range("K141").Value = 7427.000100000001
$ws.Range("L141").Value = 3976824.3
$ws.Range("M141").Value = -2247.000100000001
$ws.Range("N141").Value = -3987184.3

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32,I32,K32,M32
$ws.Range("H32").Value = 8107.758
$ws.Range("I32").Value = 7053.5
$ws.Range("K32").Value = 7053.5
$ws.Range("M32").Value = -6766.5
# Row 130: H130,J130,L130,N130
$ws.Range("H130").Value = 31000
$ws.Range("J130").Value = 31000
$ws.Range("L130").Value = 31000
$ws.Range("N130").Value = -41040
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 3687.0334
$ws.Range("I132").Value = 3354.9092
$ws.Range("K132").Value = 10064.7276
$ws.Range("M132").Value = -7534.7276

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20: H20,I20,J20,K20,L20,M20,N20
$ws.Range("H20").Value = 1897.6364
$ws.Range("I20").Value = 1853.6
$ws.Range("J20").Value = 1934.3334
$ws.Range("K20").Value = 1853.6
$ws.Range("L20").Value = 1934.3334
$ws.Range("M20").Value = -1606.6
$ws.Range("N20").Value = -2428.3334

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31,I31,J31,K31,L31,M31,N31
$ws.Range("H31").Value = 1853937.4
$ws.Range("I31").Value = 2129049.8
$ws.Range("J31").Value = 6754.143
$ws.Range("K31").Value = 2129049.8
$ws.Range("L31").Value = 6754.143
$ws.Range("M31").Value = -2128754.8
$ws.Range("N31").Value = -7344.143
# Row 34: H34,I34,J34,K34,L34,M34,N34
$ws.Range("H34").Value = 1853937.4
$ws.Range("I34").Value = 2129049.8
$ws.Range("J34").Value = 6754.143
$ws.Range("K34").Value = 2129049.8
$ws.Range("L34").Value = 6754.143
$ws.Range("M34").Value = -2128847.8
$ws.Range("N34").Value = -7158.143
# Row 86: H86,I86,J86,K86,L86,M86,N86
$ws.Range("H86").Value = 7363.2144
$ws.Range("I86").Value = 5535.75
$ws.Range("J86").Value = 9799.833000000001
$ws.Range("K86").Value = 5535.75
$ws.Range("L86").Value = 9799.833000000001
$ws.Range("M86").Value = -4412.75
$ws.Range("N86").Value = -12045.833
# Row 89: H89,I89,J89,K89,L89,M89,N89
$ws.Range("H89").Value = 7363.2144
$ws.Range("I89").Value = 5535.75
$ws.Range("J89").Value = 9799.833000000001
$ws.Range("K89").Value = 27678.75
$ws.Range("L89").Value = 48999.165
$ws.Range("M89").Value = -22062.75
$ws.Range("N89").Value = -60231.165
# Row 99: H99,I99,K99,M99
$ws.Range("H99").Value = 3040.818
$ws.Range("I99").Value = 1888.625
$ws.Range("K99").Value = 1888.625
$ws.Range("M99").Value = -390.625
# Row 107: H107,I107,K107,M107
$ws.Range("H107").Value = 1845.6
$ws.Range("I107").Value = 564.25
$ws.Range("K107").Value = 564.25
$ws.Range("M107").Value = 1355.75
# Row 126: H126,I126,K126,M126
$ws.Range("H126").Value = 3040.818
$ws.Range("I126").Value = 1888.625
$ws.Range("K126").Value = 5665.875
$ws.Range("M126").Value = -3195.875

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5: H5,I5,J5,K5,L5,M5,N5
$ws.Range("H5").Value = 859.6799999999999
$ws.Range("I5").Value = 454.18182
$ws.Range("J5").Value = 3833.3333
$ws.Range("K5").Value = 1362.54546
$ws.Range("L5").Value = 11499.9999
$ws.Range("M5").Value = -1250.54546
$ws.Range("N5").Value = -11723.9999
# Row 64: H64,J64,L64,N64
$ws.Range("H64").Value = 13275.556
$ws.Range("J64").Value = 22400
$ws.Range("L64").Value = 67200
$ws.Range("N64").Value = -67740
# Row 67: H67,J67,L67,N67
$ws.Range("H67").Value = 13275.556
$ws.Range("J67").Value = 22400
$ws.Range("L67").Value = 67200
$ws.Range("N67").Value = -69072
# Row 131: H131,I131,J131,K131,L131,M131,N131
$ws.Range("H131").Value = 1317.8209
$ws.Range("I131").Value = 3544.2856
$ws.Range("J131").Value = 1058.0667
$ws.Range("K131").Value = 10632.8568
$ws.Range("L131").Value = 3174.2001
$ws.Range("M131").Value = -5592.856800000001
$ws.Range("N131").Value = -13254.2001
# Row 135: H135,I135,J135,K135,L135,M135,N135
$ws.Range("H135").Value = 859.6799999999999
$ws.Range("I135").Value = 454.18182
$ws.Range("J135").Value = 3833.3333
$ws.Range("K135").Value = 4087.63638
$ws.Range("L135").Value = 34499.9997
$ws.Range("M135").Value = -1552.63638
$ws.Range("N135").Value = -39569.9997

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 126: H126,I126,J126,K126,L126,M126,N126
$ws.Range("H126").Value = 2824.4482
$ws.Range("I126").Value = 1293.0667
$ws.Range("J126").Value = 4465.2144
$ws.Range("K126").Value = 3879.2001
$ws.Range("L126").Value = 13395.6432
$ws.Range("M126").Value = -1409.2001
$ws.Range("N126").Value = -18335.6432
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 3990.7576
$ws.Range("I132").Value = 3934.85
$ws.Range("K132").Value = 11804.55
$ws.Range("M132").Value = -9274.549999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22: H22,I22,K22,M22
$ws.Range("H22").Value = 142859000
$ws.Range("I22").Value = 333333660
$ws.Range("K22").Value = 333333660
$ws.Range("M22").Value = -333333365
# Row 27: H27,I27,K27,M27
$ws.Range("H27").Value = 142859000
$ws.Range("I27").Value = 333333660
$ws.Range("K27").Value = 333333660
$ws.Range("M27").Value = -333333553
# Row 68: H68,I68,J68,K68,L68,M68,N68
$ws.Range("H68").Value = 1843.6364
$ws.Range("I68").Value = 1038
$ws.Range("J68").Value = 9900
$ws.Range("K68").Value = 1038
$ws.Range("L68").Value = 9900
$ws.Range("M68").Value = -289
$ws.Range("N68").Value = -11398
# Row 71: H71,I71,J71,K71,L71,M71,N71
$ws.Range("H71").Value = 1843.6364
$ws.Range("I71").Value = 1038
$ws.Range("J71").Value = 9900
$ws.Range("K71").Value = 5190
$ws.Range("L71").Value = 49500
$ws.Range("M71").Value = -1446
$ws.Range("N71").Value = -56988
# Row 128: H128,J128,L128,N128
$ws.Range("H128").Value = 27000
$ws.Range("J128").Value = 27000
$ws.Range("L128").Value = 27000
$ws.Range("N128").Value = -36960

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62: H62,I62,J62,K62,L62,M62,N62
$ws.Range("H62").Value = 4141
$ws.Range("I62").Value = 3536.6
$ws.Range("J62").Value = 4476.778
$ws.Range("K62").Value = 3536.6
$ws.Range("L62").Value = 4476.778
$ws.Range("M62").Value = -2912.6
$ws.Range("N62").Value = -5724.778
# Row 65: H65,I65,J65,K65,L65,M65
$ws.Range("H65").Value = 4141
$ws.Range("I65").Value = 3536.6
$ws.Range("J65").Value = 4476.778
$ws.Range("K65").Value = 17683
$ws.Range("L65").Value = 22383.89
$ws.Range("M65").Value = -14563
# Row 107: H107,I107,K107,M107
$ws.Range("H107").Value = 710.0417
$ws.Range("I107").Value = 301.2353
$ws.Range("K107").Value = 903.7058999999999
$ws.Range("M107").Value = 1016.2941
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 162559.88
$ws.Range("I132").Value = 185920.53
$ws.Range("K132").Value = 557761.59
$ws.Range("M132").Value = -555231.59

